$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Add the new final slide (slide 4), re-using the layout of slide 3 ("General
# Slide" - a blank layout), matching "Ajout de la derniere slide".
# ---------------------------------------------------------------------------
$s = $p.Slides.Add($p.Slides.Count + 1, 5)

# --- Shape 1: big title-ish text box ("C'est moi qui ai trouvé ISO-8859-1, voilà !") ---
$tb1 = $s.Shapes.AddTextbox(1, 131.7735433070866, 67.92456692913386, 1191.3962992125985, 201.1452755905512)
$tb1.Name = "ZoneTexte 1"
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $true

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "C" + [char]0x2019 + "est moi qui ai trouv" + [char]0x00E9 + " ISO-8859-1, voil" + [char]0x00E0 + " !"
$tr1.Font.Size = 80
$tr1.LanguageID = "fr-FR"

$tb1.TextFrame.AutoSize = 1
$tb1.Height = 201.1452755905512

# --- Shape 2: small thank-you note ("Merci <smiley>") ---
$tb2 = $s.Shapes.AddTextbox(1, 1268.830157480315, 643.9245669291339, 139.92456692913385, 29.081259842519685)
$tb2.Name = "ZoneTexte 2"
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $true

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "Merci " + [string][char]0x4A
$tr2.Font.Size = 18
$tr2.LanguageID = "fr-FR"

$smiley = $tr2.Characters(7, 1)
$smiley.Font.Name = "Wingdings"

$tb2.TextFrame.AutoSize = 1
$tb2.Height = 29.081259842519685
